$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to stay text-typed
# (mirrors the existing inlineStr storage) so values keep exact formatting
# such as trailing zeros instead of being auto-parsed into numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.372.69"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "1.687.40"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.75%  "

$ws.Range("D5").Value = "218.47"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").Value = "0.5465"
$ws.Range("E6").Value = "  +4.33%  "

$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("D8").Value = "0.2718"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("D9").Value = "0.06463"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").Value = "22.03"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").Value = "0.07690"
$ws.Range("E11").Value = "  +3.14%  "

$ws.Range("D12").Value = "1.706.67"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").Value = "4.531"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "0.5813"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").Value = "0.000008383"
$ws.Range("E15").Value = "  -2.01%  "

$ws.Range("D16").Value = "65.06"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "26.417.69"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "4.936"
$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("D20").Value = "10.98"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").Value = "6.232"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("D24").Value = "149.69"
$ws.Range("E24").Value = "  +3.04%  "

$ws.Range("D25").Value = "0.1318"
$ws.Range("E25").Value = "  +5.68%  "

$ws.Range("D26").Value = "7.873"
$ws.Range("E26").Value = "  +2.91%  "

$ws.Range("D27").Value = "15.69"
$ws.Range("E27").Value = "  -0.96%  "

$ws.Range("D28").Value = "0.06300"
$ws.Range("E28").Value = "  -6.31%  "

$ws.Range("D29").Value = "1.406"
$ws.Range("E29").Value = "  +5.27%  "

$ws.Range("D30").Value = "1.327"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "3.581"
$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").Value = "3.578"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").Value = "1.676"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("E34").Value = "  +1.50%  "

$ws.Range("D35").Value = "0.6150"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").Value = "2.413"
$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D38").Value = "6.246"
$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("D39").Value = "1.113.44"
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").Value = "0.01623"
$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("D41").Value = "0.8815"
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").Value = "101.41"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "1.839.76"
$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -2.00%  "

$ws.Range("D46").Value = "57.31"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").Value = "8.175"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "0.05272"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Value = "6.039"
$ws.Range("E51").Value = "  +0.39%  "
